# Fruta / hortaliza, semanal
# Re-applies the weekly refresh of the "Fecha" (D), "Volumen" (M) and the
# price columns (N/O/P min-max-weighted-avg, S $/Kg) for the
# "Terminal La Palmera de La Serena - Coco" sheet. Row 1 (headers) and
# row 25 are left untouched; all other data rows (2-24, 26-30) get the
# refreshed figures below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value2 = 44407
$ws.Range("M2").Value2 = 160
$ws.Range("N2").Value2 = 20000
$ws.Range("O2").Value2 = 21000
$ws.Range("P2").Value2 = 20500
$ws.Range("S2").Value2 = 1025

$ws.Range("D3").Value2 = 44365
$ws.Range("M3").Value2 = 100
$ws.Range("N3").Value2 = 20000
$ws.Range("O3").Value2 = 21000
$ws.Range("P3").Value2 = 20500
$ws.Range("S3").Value2 = 1025

$ws.Range("D4").Value2 = 44420
$ws.Range("M4").Value2 = 160
$ws.Range("N4").Value2 = 20000
$ws.Range("O4").Value2 = 21000
$ws.Range("P4").Value2 = 20500
$ws.Range("S4").Value2 = 1025

$ws.Range("D5").Value2 = 44333
$ws.Range("M5").Value2 = 100
$ws.Range("N5").Value2 = 19500
$ws.Range("O5").Value2 = 20000
$ws.Range("P5").Value2 = 19750
$ws.Range("S5").Value2 = 988

$ws.Range("D6").Value2 = 44418
$ws.Range("M6").Value2 = 200
$ws.Range("N6").Value2 = 20000
$ws.Range("O6").Value2 = 21000
$ws.Range("P6").Value2 = 20500
$ws.Range("S6").Value2 = 1025

$ws.Range("D7").Value2 = 44343
$ws.Range("M7").Value2 = 100
$ws.Range("N7").Value2 = 19500
$ws.Range("O7").Value2 = 20000
$ws.Range("P7").Value2 = 19750
$ws.Range("S7").Value2 = 988

$ws.Range("D8").Value2 = 44473
$ws.Range("M8").Value2 = 40
$ws.Range("N8").Value2 = 19500
$ws.Range("O8").Value2 = 20000
$ws.Range("P8").Value2 = 19750
$ws.Range("S8").Value2 = 988

$ws.Range("D9").Value2 = 44448
$ws.Range("M9").Value2 = 100
$ws.Range("N9").Value2 = 20000
$ws.Range("O9").Value2 = 21000
$ws.Range("P9").Value2 = 20500
$ws.Range("S9").Value2 = 1025

$ws.Range("D10").Value2 = 44301
$ws.Range("M10").Value2 = 100
$ws.Range("N10").Value2 = 18000
$ws.Range("O10").Value2 = 19000
$ws.Range("P10").Value2 = 18500
$ws.Range("S10").Value2 = 925

$ws.Range("D11").Value2 = 44427
$ws.Range("M11").Value2 = 200
$ws.Range("N11").Value2 = 20000
$ws.Range("O11").Value2 = 21000
$ws.Range("P11").Value2 = 20500
$ws.Range("S11").Value2 = 1025

$ws.Range("D12").Value2 = 44467
$ws.Range("M12").Value2 = 200
$ws.Range("N12").Value2 = 20000
$ws.Range("O12").Value2 = 21000
$ws.Range("P12").Value2 = 20500
$ws.Range("S12").Value2 = 1025

$ws.Range("D13").Value2 = 44326
$ws.Range("M13").Value2 = 160
$ws.Range("N13").Value2 = 19500
$ws.Range("O13").Value2 = 20000
$ws.Range("P13").Value2 = 19750
$ws.Range("S13").Value2 = 988

$ws.Range("D14").Value2 = 44435
$ws.Range("M14").Value2 = 260
$ws.Range("N14").Value2 = 20000
$ws.Range("O14").Value2 = 22000
$ws.Range("P14").Value2 = 21115
$ws.Range("S14").Value2 = 1056

$ws.Range("D15").Value2 = 44336
$ws.Range("M15").Value2 = 100
$ws.Range("N15").Value2 = 19500
$ws.Range("O15").Value2 = 20000
$ws.Range("P15").Value2 = 19750
$ws.Range("S15").Value2 = 988

$ws.Range("D16").Value2 = 44350
$ws.Range("M16").Value2 = 160
$ws.Range("N16").Value2 = 19000
$ws.Range("O16").Value2 = 20000
$ws.Range("P16").Value2 = 19500
$ws.Range("S16").Value2 = 975

$ws.Range("D17").Value2 = 44434
$ws.Range("M17").Value2 = 100
$ws.Range("N17").Value2 = 20000
$ws.Range("O17").Value2 = 21000
$ws.Range("P17").Value2 = 20500
$ws.Range("S17").Value2 = 1025

$ws.Range("D18").Value2 = 44417
$ws.Range("M18").Value2 = 160
$ws.Range("N18").Value2 = 20000
$ws.Range("O18").Value2 = 21000
$ws.Range("P18").Value2 = 20500
$ws.Range("S18").Value2 = 1025

$ws.Range("D19").Value2 = 44410
$ws.Range("M19").Value2 = 200
$ws.Range("N19").Value2 = 20000
$ws.Range("O19").Value2 = 21000
$ws.Range("P19").Value2 = 20500
$ws.Range("S19").Value2 = 1025

$ws.Range("D20").Value2 = 44466
$ws.Range("M20").Value2 = 100
$ws.Range("N20").Value2 = 20000
$ws.Range("O20").Value2 = 21000
$ws.Range("P20").Value2 = 20500
$ws.Range("S20").Value2 = 1025

$ws.Range("D21").Value2 = 44462
$ws.Range("M21").Value2 = 100
$ws.Range("N21").Value2 = 19500
$ws.Range("O21").Value2 = 20000
$ws.Range("P21").Value2 = 19750
$ws.Range("S21").Value2 = 988

$ws.Range("D22").Value2 = 44364
$ws.Range("M22").Value2 = 140
$ws.Range("N22").Value2 = 20000
$ws.Range("O22").Value2 = 21000
$ws.Range("P22").Value2 = 20500
$ws.Range("S22").Value2 = 1025

$ws.Range("D23").Value2 = 44431
$ws.Range("M23").Value2 = 160
$ws.Range("N23").Value2 = 21000
$ws.Range("O23").Value2 = 22000
$ws.Range("P23").Value2 = 21500
$ws.Range("S23").Value2 = 1075

$ws.Range("D24").Value2 = 44442
$ws.Range("M24").Value2 = 140
$ws.Range("N24").Value2 = 20000
$ws.Range("O24").Value2 = 21000
$ws.Range("P24").Value2 = 20500
$ws.Range("S24").Value2 = 1025

$ws.Range("D25").Value2 = 44445
$ws.Range("M25").Value2 = 160
$ws.Range("N25").Value2 = 20000
$ws.Range("O25").Value2 = 21000
$ws.Range("P25").Value2 = 20500
$ws.Range("S25").Value2 = 1025

$ws.Range("D26").Value2 = 44441
$ws.Range("M26").Value2 = 160
$ws.Range("N26").Value2 = 20000
$ws.Range("O26").Value2 = 21000
$ws.Range("P26").Value2 = 20500
$ws.Range("S26").Value2 = 1025

$ws.Range("D27").Value2 = 44335
$ws.Range("M27").Value2 = 200
$ws.Range("N27").Value2 = 19000
$ws.Range("O27").Value2 = 20000
$ws.Range("P27").Value2 = 19500
$ws.Range("S27").Value2 = 975

$ws.Range("D28").Value2 = 44315
$ws.Range("M28").Value2 = 100
$ws.Range("N28").Value2 = 20000
$ws.Range("O28").Value2 = 21000
$ws.Range("P28").Value2 = 20500
$ws.Range("S28").Value2 = 1025

$ws.Range("D29").Value2 = 44474
$ws.Range("M29").Value2 = 200
$ws.Range("N29").Value2 = 19000
$ws.Range("O29").Value2 = 20000
$ws.Range("P29").Value2 = 19500
$ws.Range("S29").Value2 = 975

$ws.Range("D30").Value2 = 44428
$ws.Range("M30").Value2 = 100
$ws.Range("N30").Value2 = 20000
$ws.Range("O30").Value2 = 21000
$ws.Range("P30").Value2 = 20500
$ws.Range("S30").Value2 = 1025
